$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegisterUser")

$ws.Range("A2").Value = "Shruti"
$ws.Range("B2").Value = "Malhotra"
$ws.Range("C2").Value = "shruti.malhotra@gmail.com"
$ws.Range("D2").Value = "shruti.malhotra@gmail.com"
